# Update the "Förändrad" (Changed) date column (C) from 45202 to 45203
# for rows 2 through 116 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C116").Value = 45203
